$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): update 想去人数 (want-to-go count) for rows 2 and 3
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5217
$ws1.Range("F3").Value = 160

# Sheet "全部类型" (All types): same two cells mirror the exhibitions sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5217
$ws4.Range("F3").Value = 160
